# Fix the malformed email for "Sai Sirisha Devineni" (row 5, column A):
# "s531367asdfnwmissouri.edu" -> "s531367@nwmissouri.edu"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $null
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($val -eq "s531367asdfnwmissouri.edu") {
        $target = $r
        break
    }
}

if ($target -ne $null) {
    $ws.Cells.Item($target, 1).Value2 = "s531367@nwmissouri.edu"
}

# Restore the view/selection state recorded by the instructor dashboard fix:
# top row scrolled to A4, active selection on B23.
$ws.Range("B23").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
